# Add new row 5 with pasta fusi product code, matching the pattern of
# existing rows (A: description, B: code, C: intentionally blank cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "25 lb pasta fusi for ravioli"
$ws.Range("B5").Value = "25 01TRAV"

# Materialize an (empty) C5 cell - mirrors C3/C4 which are present but blank -
# without assigning any explicit (non-default) style.
$ws.Range("C5").Style = "Normal"
